$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Join the two runs that were split by the stray "_GoBack" bookmark in the
#    1Tim 5:19 quotation, and drop that bookmark from this location.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Do not receive an accusation aga*inst an elder except from two or three witnesses.",
    $false, $false, $true, $false, $false, $true, 1, $false,
    "Do not receive an accusation against an elder except from two or three witnesses.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Fill in the worksheet's empty "Question" paragraphs with their text.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(9).Range.Text =
    "Why is it important for Paul to back up his statement in 5:17 with citations from the Old Testament in 5:18?"

$d.Paragraphs.Item(13).Range.Text =
    "What is the purpose for requiring two or three witnesses to corroborate accusations of an elder in 5:19?"

$d.Paragraphs.Item(17).Range.Text =
    "Who are the " + [char]0x201C + "all" + [char]0x201D + " before whom sinning elders are to be rebuked, and who are the " + [char]0x201C + "rest" + [char]0x201D + " who should respond in fear (5:20)?"

$d.Paragraphs.Item(21).Range.Text =
    "Why is it so important for Timothy to act without prejudice or partiality in judging accusations against elders (5:21)?"

$d.Paragraphs.Item(25).Range.Text =
    "What does it mean to " + [char]0x201C + "not lay hands on anyone hastily," + [char]0x201D + " and why was that important (5:22)?"

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark: it now marks the last blank "Lines"
#    paragraph right before the 5:24-25 question.
# ---------------------------------------------------------------------------
$goBackRange = $d.Paragraphs.Item(28).Range
$goBackRange.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$d.Paragraphs.Item(29).Range.Text =
    "In what way does the observation in 5:24" + [char]0x2013 + "25 support the injunction in 5:22?"

$d.Paragraphs.Item(33).Range.Text =
    "Why might a Christian slave be tempted to insubordination to his human master? How would that cause " + [char]0x201C + "the name of God and His doctrine" + [char]0x201D + " to be " + [char]0x201C + "blasphemed" + [char]0x201D + " (6:1)?"

$d.Paragraphs.Item(37).Range.Text =
    "Why does Paul tolerate the institution of slavery (6:1" + [char]0x2013 + "2) rather than calling Timothy to preach a political message of liberation?"
